# Correzione di alcune exit condition
# The "exit condition" paragraph currently reads:
#   "L'utente non riesce a comunicare col server e non riesce a confermare l'ordine."
# and must become:
#   "L'utente non riesce a confermare l'ordine."
# Only the text of the first run changes; the other two runs (and their
# rPr/rsid metadata) must stay exactly as they are, so we replace the whole
# paragraph via InsertXML with the precise run layout instead of a plain
# Find/Replace (which would coalesce the runs).

$d = $word.ActiveDocument

$targetParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*comunicare col server*non riesce a confermare*") {
        $targetParagraph = $candidate
        break
    }
}

if ($targetParagraph -ne $null) {
    $rng = $targetParagraph.Range

    $paraXml = '<w:p w14:paraId="3C4A4D96" w14:textId="0F592028" w:rsidR="00E12504" w:rsidRPr="00506DC0" w:rsidRDefault="00D95CCD" w:rsidP="006F6523"><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve">L’utente non riesce a </w:t></w:r><w:r w:rsidR="00C80A09"><w:rPr><w:bCs/></w:rPr><w:t>confermare l’ordine</w:t></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:t>.</w:t></w:r></w:p>'

    $packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $paraXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($packageXml)
}
